$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (DAMSLTag, DialogAct)
$updates = @{
    3  = @("%", "Uninterpretable")
    10 = @("sd", "Statement-non-opinion")
    13 = @("ba", "Appreciation")
    19 = @("ba", "Appreciation")
    21 = @("sd", "Statement-non-opinion")
    25 = @("sd", "Statement-non-opinion")
    45 = @("%", "Uninterpretable")
    46 = @("sv", "Statement-opinion")
    54 = @("sd", "Statement-non-opinion")
    58 = @("sv", "Statement-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}

$wb.Save()
